$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.084.33"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").Value = "3.569.66"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'241.00"
$ws.Range("E5").Value = "  +2.34%  "

$ws.Range("D6").Value = "'654.31"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'1.72"
$ws.Range("E7").Value = "  +16.33%  "

$ws.Range("D8").Value = "'0.424"
$ws.Range("E8").Value = "  +6.32%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").Value = "'1.05"
$ws.Range("E10").Value = "  +4.56%  "

$ws.Range("D11").Value = "3.566.86"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").Value = "'44.05"
$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("D13").Value = "'0.203"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").Value = "'6.42"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "4.231.46"
$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("D16").Value = "96.802.26"
$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("E17").Value = "  +2.45%  "

$ws.Range("D18").Value = "'8.63"
$ws.Range("E18").Value = "  +11.43%  "

$ws.Range("D19").Value = "3.564.31"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "'12.69"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").Value = "'17.93"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").Value = "'0.530"
$ws.Range("E22").Value = "  +11.02%  "

$ws.Range("D23").Value = "'3.47"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").Value = "'512.69"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").Value = "'0.0000205"
$ws.Range("E25").Value = "  +5.12%  "

$ws.Range("D26").Value = "'6.82"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").Value = "'101.33"
$ws.Range("E27").Value = "  +6.39%  "

$ws.Range("D28").Value = "'13.00"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("D29").Value = "3.758.41"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "'0.167"
$ws.Range("E30").Value = "  +16.59%  "

$ws.Range("D31").Value = "'2.99"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").Value = "'11.85"
$ws.Range("E32").Value = "  +2.92%  "

$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("D34").Value = "'0.183"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").Value = "'31.69"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").Value = "'8.80"
$ws.Range("E37").Value = "  +2.87%  "

$ws.Range("D38").Value = "'614.46"
$ws.Range("E38").Value = "  +5.20%  "

$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("E40").Value = "  -3.39%  "

$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'1.95"
$ws.Range("E41").Value = "  +7.21%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.154"
$ws.Range("E42").Value = "  +2.33%  "

$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").Value = "'0.918"
$ws.Range("E44").Value = "  +1.66%  "

$ws.Range("D45").Value = "'6.00"
$ws.Range("E45").Value = "  +4.51%  "

$ws.Range("D46").Value = "'0.0437"
$ws.Range("E46").Value = "  +5.25%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").Value = "'0.407"
$ws.Range("E49").Value = "  +31.28%  "

$ws.Range("D50").Value = "'8.51"
$ws.Range("E50").Value = "  +4.14%  "

$ws.Range("D51").Value = "'32.99"
$ws.Range("E51").Value = "  -2.78%  "
